$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 2
    4  = 1
    5  = 1
    6  = 3
    7  = 1
    8  = 1
    9  = 0
    10 = 2
    11 = 0
    12 = 0
    13 = 0
    14 = 1
    15 = 0
    16 = 1
    18 = 1
    19 = 1
    20 = 0
    21 = 0
    22 = 2
    23 = 0
    24 = 2
    26 = 1
    27 = 1
    28 = 0
    29 = 1
    30 = 0
    31 = 2
    32 = 3
    33 = 1
    34 = 1
    35 = 1
    37 = 0
    38 = 0
    39 = 0
    41 = 0
    42 = 2
    43 = 0
    44 = 2
    45 = 0
    46 = 0
    47 = 0
    48 = 1
    49 = 2
    50 = 0
    51 = 0
    52 = 0
    53 = 1
    54 = 2
    55 = 1
    56 = 0
    57 = 3
    58 = 1
    59 = 0
    60 = 0
    61 = 2
    62 = 1
    63 = 0
    64 = 1
    65 = 0
    66 = 0
    67 = 0
    68 = 0
    69 = 4
    70 = 0
    71 = 0
    72 = 0
    73 = 1
    74 = 0
    75 = 2
    76 = 1
    77 = 1
    78 = 2
    79 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
